$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "discretiser_type"
$ws.Range("B15").Value = "equalfixed"

$ws.Range("A16").Value = "ohe"
$ws.Range("B16").Value = "grade"

$ws.Range("A17").Value = "outlier"
$ws.Range("B17").Value = "installment"

$ws.Range("A18").Value = "scale"
$ws.Range("B18").Value = "minmax"

$ws.Range("A19").Select()
